$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the F-column "time_taken" timestamps on the "data" sheet.
# ---------------------------------------------------------------------------
$timestamps = @(
    "2021-10-05 14:34:33.934743",
    "2021-10-05 14:34:33.934768",
    "2021-10-05 14:34:33.934772",
    "2021-10-05 14:34:33.934809",
    "2021-10-05 14:34:33.934814",
    "2021-10-05 14:34:33.934817",
    "2021-10-05 14:34:33.934820",
    "2021-10-05 14:34:33.934822",
    "2021-10-05 14:34:33.934825",
    "2021-10-05 14:34:33.934828",
    "2021-10-05 14:34:33.934831",
    "2021-10-05 14:34:33.934834",
    "2021-10-05 14:34:33.934836",
    "2021-10-05 14:34:33.934839",
    "2021-10-05 14:34:33.934842",
    "2021-10-05 14:34:33.934845",
    "2021-10-05 14:34:33.934848",
    "2021-10-05 14:34:33.934851",
    "2021-10-05 14:34:33.934853",
    "2021-10-05 14:34:33.934856",
    "2021-10-05 14:34:33.934859",
    "2021-10-05 14:34:33.934862",
    "2021-10-05 14:34:33.934890",
    "2021-10-05 14:34:33.934924",
    "2021-10-05 14:34:33.934950",
    "2021-10-05 14:34:33.934956",
    "2021-10-05 14:34:33.934961",
    "2021-10-05 14:34:33.934970",
    "2021-10-05 14:34:33.934975",
    "2021-10-05 14:34:33.934979",
    "2021-10-05 14:34:33.934983"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" worksheet after "data".
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the bold/bordered header style used on the "data" sheet (style index 1)
# onto the header row and the A2 index cell of the new sheet.
$dataSheet.Range("B1").Copy($metaSheet.Range("B1:G1"))
$dataSheet.Range("A2").Copy($metaSheet.Range("A2"))

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Mandibulofacial Acrofacial dysostosis"
$metaSheet.Range("C2").Value = 136
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.0"
$metaSheet.Range("E2").Value = "2021-08-08T08:27:23.150341Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:34:33.928772"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/136/?format=json"

$metaSheet.Range("A1").Select()
